$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The "date" and "time" field types are no longer used by this form -
# every date/time field becomes a plain "text" field (JGI app change).
$ws.Cells.Item(2, 3).Value = "text"
$ws.Cells.Item(4, 3).Value = "text"
$ws.Cells.Item(6, 3).Value = "text"
